# Update the "想去人数" (F) and occasionally "最低票价" (G) figures on the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets to reflect newly scraped
# totals, as captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 166
$ws1.Range("F4").Value  = 169
$ws1.Range("F5").Value  = 4812
$ws1.Range("F7").Value  = 42
$ws1.Range("F8").Value  = 533
$ws1.Range("F9").Value  = 483
$ws1.Range("F10").Value = 27
$ws1.Range("F12").Value = 1356
$ws1.Range("F13").Value = 3107
$ws1.Range("F14").Value = 392
$ws1.Range("F15").Value = 115
$ws1.Range("G15").Value = 58
$ws1.Range("F16").Value = 95
$ws1.Range("F18").Value = 2462
$ws1.Range("F19").Value = 114
$ws1.Range("F23").Value = 22
$ws1.Range("F24").Value = 116

# --- 全部类型 sheet ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 166
$ws4.Range("F4").Value  = 169
$ws4.Range("F6").Value  = 4812
$ws4.Range("F8").Value  = 42
$ws4.Range("F9").Value  = 533
$ws4.Range("F10").Value = 483
$ws4.Range("F11").Value = 27
$ws4.Range("F13").Value = 1356
$ws4.Range("F14").Value = 3107
$ws4.Range("F15").Value = 392
$ws4.Range("F16").Value = 115
$ws4.Range("G16").Value = 58
$ws4.Range("F17").Value = 95
$ws4.Range("F19").Value = 2462
$ws4.Range("F20").Value = 114
$ws4.Range("F24").Value = 22
$ws4.Range("F25").Value = 116
